# Applies "Atualizacao de bases das ligas" update (2024-05-02 20:28)
# to "Germany Bundesliga I" sheet: refreshed match ids/teams/odds for
# rows 21-24, 39-42, 281, 283-285, 288-289.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Cells.Item(21, 2).Value = 6847566
$ws.Cells.Item(21, 5).Value = "Werder Bremen"
$ws.Cells.Item(21, 6).Value = "Mainz"
$ws.Cells.Item(21, 7).Value = 4
$ws.Cells.Item(21, 10).Value = 2.875
$ws.Cells.Item(21, 11).Value = 3.6
$ws.Cells.Item(21, 12).Value = 2.3
$ws.Cells.Item(21, 13).Value = 2.6
$ws.Cells.Item(21, 14).Value = 3.4
$ws.Cells.Item(21, 15).Value = 2.6
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 1.95
$ws.Cells.Item(21, 18).Value = 1.95
$ws.Cells.Item(21, 20).Value = 2.04
$ws.Cells.Item(21, 21).Value = 1.86
$ws.Cells.Item(21, 22).Value = 1.6
$ws.Cells.Item(21, 25).Value = 0.95
$ws.Cells.Item(21, 27).Value = 1.04

# Row 22
$ws.Cells.Item(22, 2).Value = 6847565
$ws.Cells.Item(22, 5).Value = "TSG Hoffenheim"
$ws.Cells.Item(22, 6).Value = "Wolfsburg"
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 10).Value = 2.3
$ws.Cells.Item(22, 11).Value = 3.6
$ws.Cells.Item(22, 12).Value = 2.8
$ws.Cells.Item(22, 13).Value = 2.5
$ws.Cells.Item(22, 14).Value = 3.75
$ws.Cells.Item(22, 15).Value = 2.55
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 1.95
$ws.Cells.Item(22, 18).Value = 1.95
$ws.Cells.Item(22, 20).Value = 1.95
$ws.Cells.Item(22, 21).Value = 1.95
$ws.Cells.Item(22, 22).Value = 1.5
$ws.Cells.Item(22, 25).Value = 0.95
$ws.Cells.Item(22, 27).Value = 0.95

# Row 23
$ws.Cells.Item(23, 2).Value = 6847563
$ws.Cells.Item(23, 5).Value = "Bayer Leverkusen"
$ws.Cells.Item(23, 6).Value = "Darmstadt"
$ws.Cells.Item(23, 7).Value = 5
$ws.Cells.Item(23, 10).Value = 1.4
$ws.Cells.Item(23, 11).Value = 4.75
$ws.Cells.Item(23, 12).Value = 7.5
$ws.Cells.Item(23, 13).Value = 1.222
$ws.Cells.Item(23, 14).Value = 6.5
$ws.Cells.Item(23, 15).Value = 11
$ws.Cells.Item(23, 16).Value = -2
$ws.Cells.Item(23, 17).Value = 2.07
$ws.Cells.Item(23, 18).Value = 1.83
$ws.Cells.Item(23, 20).Value = 1.99
$ws.Cells.Item(23, 21).Value = 1.91
$ws.Cells.Item(23, 22).Value = 0.222
$ws.Cells.Item(23, 25).Value = 1.07
$ws.Cells.Item(23, 27).Value = 0.99

# Row 24
$ws.Cells.Item(24, 2).Value = 6847568
$ws.Cells.Item(24, 5).Value = "VfB Stuttgart"
$ws.Cells.Item(24, 6).Value = "SC Freiburg"
$ws.Cells.Item(24, 7).Value = 5
$ws.Cells.Item(24, 10).Value = 2
$ws.Cells.Item(24, 11).Value = 3.4
$ws.Cells.Item(24, 12).Value = 3.8
$ws.Cells.Item(24, 13).Value = 2.25
$ws.Cells.Item(24, 14).Value = 3.5
$ws.Cells.Item(24, 15).Value = 3.1
$ws.Cells.Item(24, 16).Value = -0.25
$ws.Cells.Item(24, 17).Value = 2
$ws.Cells.Item(24, 18).Value = 1.9
$ws.Cells.Item(24, 20).Value = 1.98
$ws.Cells.Item(24, 21).Value = 1.92
$ws.Cells.Item(24, 22).Value = 1.25
$ws.Cells.Item(24, 25).Value = 1
$ws.Cells.Item(24, 27).Value = 0.98

# Row 39
$ws.Cells.Item(39, 2).Value = 6847578
$ws.Cells.Item(39, 5).Value = "Borussia Dortmund"
$ws.Cells.Item(39, 6).Value = "Wolfsburg"
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = "H"
$ws.Cells.Item(39, 10).Value = 1.615
$ws.Cells.Item(39, 11).Value = 4
$ws.Cells.Item(39, 12).Value = 5
$ws.Cells.Item(39, 13).Value = 1.727
$ws.Cells.Item(39, 14).Value = 4.2
$ws.Cells.Item(39, 15).Value = 4.333
$ws.Cells.Item(39, 16).Value = -0.75
$ws.Cells.Item(39, 17).Value = 1.92
$ws.Cells.Item(39, 18).Value = 1.98
$ws.Cells.Item(39, 19).Value = 3.25
$ws.Cells.Item(39, 20).Value = 1.99
$ws.Cells.Item(39, 21).Value = 1.91
$ws.Cells.Item(39, 22).Value = 0.7270000000000001
$ws.Cells.Item(39, 24).Value = -1
$ws.Cells.Item(39, 25).Value = 0.46
$ws.Cells.Item(39, 26).Value = -0.5
$ws.Cells.Item(39, 28).Value = 0.9099999999999999

# Row 40
$ws.Cells.Item(40, 2).Value = 6847579
$ws.Cells.Item(40, 5).Value = "Union Berlin"
$ws.Cells.Item(40, 6).Value = "TSG Hoffenheim"
$ws.Cells.Item(40, 8).Value = 2
$ws.Cells.Item(40, 10).Value = 1.909
$ws.Cells.Item(40, 11).Value = 3.75
$ws.Cells.Item(40, 12).Value = 3.6
$ws.Cells.Item(40, 13).Value = 2.05
$ws.Cells.Item(40, 14).Value = 3.6
$ws.Cells.Item(40, 15).Value = 3.3
$ws.Cells.Item(40, 16).Value = -0.25
$ws.Cells.Item(40, 17).Value = 1.83
$ws.Cells.Item(40, 18).Value = 2.07
$ws.Cells.Item(40, 19).Value = 2.75
$ws.Cells.Item(40, 20).Value = 1.97
$ws.Cells.Item(40, 21).Value = 1.93
$ws.Cells.Item(40, 24).Value = 2.3
$ws.Cells.Item(40, 25).Value = -1
$ws.Cells.Item(40, 26).Value = 1.07
$ws.Cells.Item(40, 28).Value = 0.9299999999999999

# Row 41
$ws.Cells.Item(41, 2).Value = 6847581
$ws.Cells.Item(41, 5).Value = "Borussia Mgladbach"
$ws.Cells.Item(41, 6).Value = "RB Leipzig"
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 9).Value = "A"
$ws.Cells.Item(41, 10).Value = 4
$ws.Cells.Item(41, 11).Value = 4.2
$ws.Cells.Item(41, 12).Value = 1.75
$ws.Cells.Item(41, 13).Value = 4.333
$ws.Cells.Item(41, 14).Value = 4.333
$ws.Cells.Item(41, 15).Value = 1.65
$ws.Cells.Item(41, 16).Value = 0.75
$ws.Cells.Item(41, 17).Value = 2.05
$ws.Cells.Item(41, 18).Value = 1.85
$ws.Cells.Item(41, 19).Value = 3.25
$ws.Cells.Item(41, 20).Value = 2.02
$ws.Cells.Item(41, 21).Value = 1.88
$ws.Cells.Item(41, 22).Value = -1
$ws.Cells.Item(41, 24).Value = 0.6499999999999999
$ws.Cells.Item(41, 25).Value = -0.5
$ws.Cells.Item(41, 26).Value = 0.425
$ws.Cells.Item(41, 27).Value = -1
$ws.Cells.Item(41, 28).Value = 0.8799999999999999

# Row 42
$ws.Cells.Item(42, 2).Value = 6847583
$ws.Cells.Item(42, 5).Value = "Augsburg"
$ws.Cells.Item(42, 6).Value = "Mainz"
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 1
$ws.Cells.Item(42, 10).Value = 2.625
$ws.Cells.Item(42, 11).Value = 3.5
$ws.Cells.Item(42, 12).Value = 2.5
$ws.Cells.Item(42, 13).Value = 2.5
$ws.Cells.Item(42, 14).Value = 3.4
$ws.Cells.Item(42, 15).Value = 2.75
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(42, 17).Value = 1.82
$ws.Cells.Item(42, 18).Value = 2.08
$ws.Cells.Item(42, 19).Value = 2.5
$ws.Cells.Item(42, 20).Value = 1.89
$ws.Cells.Item(42, 21).Value = 2.01
$ws.Cells.Item(42, 22).Value = 1.5
$ws.Cells.Item(42, 25).Value = 0.8200000000000001
$ws.Cells.Item(42, 26).Value = -1
$ws.Cells.Item(42, 27).Value = 0.8899999999999999
$ws.Cells.Item(42, 28).Value = -1

# Row 281
$ws.Cells.Item(281, 13).Value = 5.25
$ws.Cells.Item(281, 14).Value = 5
$ws.Cells.Item(281, 15).Value = 1.5
$ws.Cells.Item(281, 17).Value = 1.92
$ws.Cells.Item(281, 18).Value = 1.98
$ws.Cells.Item(281, 20).Value = 1.87
$ws.Cells.Item(281, 21).Value = 2.03

# Row 283
$ws.Cells.Item(283, 2).Value = 7128365
$ws.Cells.Item(283, 5).Value = "Borussia Dortmund"
$ws.Cells.Item(283, 6).Value = "Augsburg"
$ws.Cells.Item(283, 10).Value = 1.444
$ws.Cells.Item(283, 11).Value = 4.5
$ws.Cells.Item(283, 12).Value = 6.5
$ws.Cells.Item(283, 13).Value = 1.85
$ws.Cells.Item(283, 14).Value = 4.333
$ws.Cells.Item(283, 15).Value = 3.6
$ws.Cells.Item(283, 16).Value = -0.5
$ws.Cells.Item(283, 17).Value = 1.88
$ws.Cells.Item(283, 18).Value = 2.02
$ws.Cells.Item(283, 19).Value = 3.5
$ws.Cells.Item(283, 20).Value = 1.9
$ws.Cells.Item(283, 21).Value = 2

# Row 284
$ws.Cells.Item(284, 2).Value = 7128369
$ws.Cells.Item(284, 5).Value = "Werder Bremen"
$ws.Cells.Item(284, 6).Value = "Borussia Mgladbach"
$ws.Cells.Item(284, 10).Value = 2.4
$ws.Cells.Item(284, 11).Value = 3.5
$ws.Cells.Item(284, 12).Value = 2.75
$ws.Cells.Item(284, 13).Value = 2.3
$ws.Cells.Item(284, 14).Value = 3.75
$ws.Cells.Item(284, 15).Value = 2.9
$ws.Cells.Item(284, 17).Value = 2.02
$ws.Cells.Item(284, 18).Value = 1.88
$ws.Cells.Item(284, 19).Value = 3
$ws.Cells.Item(284, 20).Value = 1.95
$ws.Cells.Item(284, 21).Value = 1.95

# Row 285
$ws.Cells.Item(285, 2).Value = 7124123
$ws.Cells.Item(285, 5).Value = "VfB Stuttgart"
$ws.Cells.Item(285, 6).Value = "Bayern Munich"
$ws.Cells.Item(285, 10).Value = 2.7
$ws.Cells.Item(285, 11).Value = 3.6
$ws.Cells.Item(285, 12).Value = 2.375
$ws.Cells.Item(285, 13).Value = 2.1
$ws.Cells.Item(285, 14).Value = 4
$ws.Cells.Item(285, 15).Value = 3.1
$ws.Cells.Item(285, 16).Value = -0.25
$ws.Cells.Item(285, 17).Value = 1.87
$ws.Cells.Item(285, 18).Value = 2.03
$ws.Cells.Item(285, 19).Value = 3.25
$ws.Cells.Item(285, 20).Value = 1.87
$ws.Cells.Item(285, 21).Value = 2.03

# Row 288
$ws.Cells.Item(288, 13).Value = 3.1
$ws.Cells.Item(288, 15).Value = 2.25
$ws.Cells.Item(288, 17).Value = 1.91
$ws.Cells.Item(288, 18).Value = 1.99
$ws.Cells.Item(288, 20).Value = 1.92
$ws.Cells.Item(288, 21).Value = 1.98

# Row 289
$ws.Cells.Item(289, 17).Value = 1.87
$ws.Cells.Item(289, 18).Value = 2.03
$ws.Cells.Item(289, 20).Value = 1.88
$ws.Cells.Item(289, 21).Value = 2.02
